$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 155
$ws.Range("I2").Value = 155
$ws.Range("K2").Value = 155
$ws.Range("M2").Value = -42

$ws.Range("H51").Value = 9123.416999999999
$ws.Range("I51").Value = 6250.5
$ws.Range("K51").Value = 6250.5
$ws.Range("M51").Value = -5766.5

$ws.Range("H53").Value = 563.9048
$ws.Range("J53").Value = 448.16666
$ws.Range("L53").Value = 448.16666
$ws.Range("N53").Value = -1722.16666

$ws.Range("H100").Value = 2479.75
$ws.Range("I100").Value = 1973.25
$ws.Range("K100").Value = 1973.25
$ws.Range("M100").Value = -1432.25

$ws.Range("H106").Value = 47621188
$ws.Range("I106").Value = 55557720
$ws.Range("K106").Value = 55557720
$ws.Range("M106").Value = -55557089

$ws.Range("H131").Value = 2145.5715
$ws.Range("I131").Value = 836.5
$ws.Range("J131").Value = 10000
$ws.Range("K131").Value = 2509.5
$ws.Range("L131").Value = 30000
$ws.Range("M131").Value = 2530.5
$ws.Range("N131").Value = -40080

$ws.Range("H138").Value = 6326.405
$ws.Range("J138").Value = 7306.543
$ws.Range("L138").Value = 21919.629
$ws.Range("N138").Value = -32199.629


# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H106").Value = 41000
$ws.Range("J106").Value = 41000
$ws.Range("L106").Value = 41000
$ws.Range("N106").Value = -43524

$ws.Range("H122").Value = 5288.528
$ws.Range("I122").Value = 4730.8696
$ws.Range("J122").Value = 6275.154
$ws.Range("K122").Value = 14192.6088
$ws.Range("L122").Value = 18825.462
$ws.Range("M122").Value = -11742.6088
$ws.Range("N122").Value = -23725.462

$ws.Range("H132").Value = 21557.258
$ws.Range("I132").Value = 31294.264
$ws.Range("K132").Value = 93882.792
$ws.Range("M132").Value = -91352.792


# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3292.4119
$ws.Range("I134").Value = 2503.5
$ws.Range("J134").Value = 5185.8
$ws.Range("K134").Value = 7510.5
$ws.Range("L134").Value = 15557.4
$ws.Range("M134").Value = -4975.5
$ws.Range("N134").Value = -20627.4

$ws.Range("H140").Value = 78233.42999999999
$ws.Range("J140").Value = 78233.42999999999
$ws.Range("L140").Value = 78233.42999999999
$ws.Range("N140").Value = -88593.42999999999


# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 11585.5
$ws.Range("J99").Value = 13100.2
$ws.Range("L99").Value = 13100.2
$ws.Range("N99").Value = -16096.2

$ws.Range("H106").Value = 26256
$ws.Range("I106").Value = 13610
$ws.Range("J106").Value = 29417.5
$ws.Range("K106").Value = 13610
$ws.Range("L106").Value = 29417.5
$ws.Range("M106").Value = -12348
$ws.Range("N106").Value = -31941.5

$ws.Range("H107").Value = 1399729.4
$ws.Range("I107").Value = 4546317
$ws.Range("J107").Value = 1246.1111
$ws.Range("K107").Value = 4546317
$ws.Range("L107").Value = 1246.1111
$ws.Range("M107").Value = -4544397
$ws.Range("N107").Value = -5086.1111

$ws.Range("H126").Value = 11585.5
$ws.Range("J126").Value = 13100.2
$ws.Range("L126").Value = 39300.60000000001
$ws.Range("N126").Value = -44240.60000000001


# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1003.7273
$ws.Range("I5").Value = 425
$ws.Range("J5").Value = 1132.3334
$ws.Range("K5").Value = 1275
$ws.Range("L5").Value = 3397.0002
$ws.Range("M5").Value = -1163
$ws.Range("N5").Value = -3621.0002

$ws.Range("H93").Value = 6117.7144
$ws.Range("I93").Value = 824
$ws.Range("J93").Value = 7000
$ws.Range("K93").Value = 2472
$ws.Range("L93").Value = 21000
$ws.Range("M93").Value = -600
$ws.Range("N93").Value = -24744

$ws.Range("H99").Value = 3302.5
$ws.Range("I99").Value = 1005
$ws.Range("J99").Value = 5600
$ws.Range("K99").Value = 3015
$ws.Range("L99").Value = 16800
$ws.Range("M99").Value = -769
$ws.Range("N99").Value = -21292

$ws.Range("H114").Value = 1898.7778
$ws.Range("I114").Value = 833
$ws.Range("J114").Value = 2111.9333
$ws.Range("K114").Value = 2499
$ws.Range("L114").Value = 6335.7999
$ws.Range("M114").Value = 755
$ws.Range("N114").Value = -12843.7999

$ws.Range("H129").Value = 1408.44
$ws.Range("I129").Value = 1008.05554
$ws.Range("K129").Value = 3024.16662
$ws.Range("M129").Value = 1975.83338

$ws.Range("H135").Value = 1003.7273
$ws.Range("I135").Value = 425
$ws.Range("J135").Value = 1132.3334
$ws.Range("K135").Value = 3825
$ws.Range("L135").Value = 10191.0006
$ws.Range("M135").Value = -1290
$ws.Range("N135").Value = -15261.0006

$ws.Range("H139").Value = 3814.2856
$ws.Range("I139").Value = 850
$ws.Range("K139").Value = 2550
$ws.Range("M139").Value = 2590


# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 33371
$ws.Range("J123").Value = 33371
$ws.Range("L123").Value = 33371
$ws.Range("N123").Value = -38271

$ws.Range("H132").Value = 4625.846
$ws.Range("I132").Value = 3931.6365
$ws.Range("J132").Value = 8444
$ws.Range("K132").Value = 11794.9095
$ws.Range("L132").Value = 25332
$ws.Range("M132").Value = -9264.9095
$ws.Range("N132").Value = -30392

$ws.Range("H137").Value = 99779
$ws.Range("J137").Value = 99779
$ws.Range("L137").Value = 99779
$ws.Range("N137").Value = -109979


# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4092.9734
$ws.Range("I7").Value = 3438.3
$ws.Range("J7").Value = 5402.32
$ws.Range("K7").Value = 3438.3
$ws.Range("L7").Value = 5402.32
$ws.Range("M7").Value = -3326.3
$ws.Range("N7").Value = -5626.32

$ws.Range("H22").Value = 1057.4615
$ws.Range("I22").Value = 899.7143
$ws.Range("K22").Value = 899.7143
$ws.Range("M22").Value = -604.7143

$ws.Range("H27").Value = 1057.4615
$ws.Range("I27").Value = 899.7143
$ws.Range("K27").Value = 899.7143
$ws.Range("M27").Value = -792.7143

$ws.Range("H46").Value = 5055.636
$ws.Range("I46").Value = 4204.25
$ws.Range("J46").Value = 5542.143
$ws.Range("K46").Value = 4204.25
$ws.Range("L46").Value = 5542.143
$ws.Range("M46").Value = -4016.25
$ws.Range("N46").Value = -5918.143

$ws.Range("H68").Value = 692388.4
$ws.Range("I68").Value = 844704.9399999999
$ws.Range("J68").Value = 6964
$ws.Range("K68").Value = 844704.9399999999
$ws.Range("L68").Value = 6964
$ws.Range("M68").Value = -843955.9399999999
$ws.Range("N68").Value = -8462

$ws.Range("H71").Value = 692388.4
$ws.Range("I71").Value = 844704.9399999999
$ws.Range("J71").Value = 6964
$ws.Range("K71").Value = 4223524.699999999
$ws.Range("L71").Value = 34820
$ws.Range("M71").Value = -4219780.699999999
$ws.Range("N71").Value = -42308

$ws.Range("H122").Value = 13040.818
$ws.Range("I122").Value = 4374.75
$ws.Range("J122").Value = 17992.857
$ws.Range("K122").Value = 13124.25
$ws.Range("L122").Value = 53978.571
$ws.Range("M122").Value = -10674.25
$ws.Range("N122").Value = -58878.571

$ws.Range("H126").Value = 4092.9734
$ws.Range("I126").Value = 3438.3
$ws.Range("J126").Value = 5402.32
$ws.Range("K126").Value = 10314.9
$ws.Range("L126").Value = 16206.96
$ws.Range("M126").Value = -7844.900000000001
$ws.Range("N126").Value = -21146.96


# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 9342.75
$ws.Range("I62").Value = 750
$ws.Range("J62").Value = 10570.286
$ws.Range("K62").Value = 750
$ws.Range("L62").Value = 10570.286
$ws.Range("M62").Value = -126
$ws.Range("N62").Value = -11818.286

$ws.Range("H65").Value = 9342.75
$ws.Range("I65").Value = 750
$ws.Range("J65").Value = 10570.286
$ws.Range("K65").Value = 3750
$ws.Range("L65").Value = 52851.43
$ws.Range("M65").Value = -630
$ws.Range("N65").Value = -59091.43

$ws.Range("H140").Value = 85056.5
$ws.Range("J140").Value = 85056.5
$ws.Range("L140").Value = 85056.5
$ws.Range("N140").Value = -95416.5

